$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped from 45205 (2023-10-06)
# to 45206 (2023-10-07) for every data row (rows 2 through 110).
$ws.Range("C2:C110").Value = 45206
